# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet a new (blank) column is inserted
# immediately before the "Late" column. Everything that used to sit in
# columns N/O/P (heading/"Outstanding" block) shifts one column to the
# right, into O/P/Q, and the freshly inserted column N is left blank
# but resized to match the width of column M next to it. The active
# sheet/selection is then moved onto "Repayment schedule" at J15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at position N (14); the old N/O/P columns
# (and all their data/styles) shift right to O/P/Q automatically.
$ws.Columns.Item(14).Insert()

# Resize the newly inserted, now-blank column N to match column M's
# width instead of leaving it at the old "Late" column's narrower
# bestFit width.
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Activate "Repayment schedule" and select J15 - this also flips
# workbookView.activeTab onto this sheet and clears tabSelected on
# whichever sheet previously had it ("Transactions").
$ws.Activate()
$ws.Range("J15").Select()
